$d = $word.ActiveDocument

# The document contains several "<id>...</id>" markers, each originally
# split across three runs: "<id>" (Courier New), the inner id text
# (Arial), and "</id>" (Courier New). This commit collapses each of the
# two "tl"/"tc"/"tcn" id markers (p122r_1 and p122r_2 - NOT the
# "fig_p122r_*" ones) into a single run containing the full literal text,
# using the formatting of the opening "<id>" run (Courier New, color
# 7f6000, sz/szCs 18).
#
# Word's Find & Replace naturally merges a multi-run match into a single
# run that carries the formatting of the first character of the match,
# which reproduces exactly that end state.

$targets = @("p122r_1", "p122r_2")

foreach ($t in $targets) {
    $old = "<id>" + $t + "</id>"
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $old
    $find.Replacement.Text = $old
    $find.Forward = $true
    $find.Wrap = 1
    $find.Format = $false
    $find.MatchCase = $true
    $find.MatchWholeWord = $false
    $find.MatchWildcards = $false
    $find.MatchSoundsLike = $false
    $find.MatchAllWordForms = $false
    [void]$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)
}
